# edit.ps1
#
# Reproduces:
#   1. The table on slide 5 switching its table style from the custom
#      "Table_0" style ({79E1BE0B-3376-428E-84A9-6EF21D299902}) to the
#      built-in style {5C0394F6-BC31-4A9A-AFA3-90D6079E29B6}.
#   2. The presentation's theme colors flipping from the "Integral"
#      (Red Violet) palette to the default "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style change (slide 5, shape 2 is the table) -----------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tbl = $tableShape.Table
$tbl.ApplyStyle("{5C0394F6-BC31-4A9A-AFA3-90D6079E29B6}")

# --- 2. Theme color swap (Integral/Red Violet -> Office Theme) -------------
function ToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Office Theme palette, in dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order
$officeThemeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47",
    "0563C1", "954F72"
)

$themeColorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColorScheme.Count; $i++) {
    $themeColorScheme.Colors($i).RGB = ToComRgb $officeThemeColors[$i - 1]
}
